$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the replacement data block for rows 45-110 (A:T)
$arr = New-Object 'object[,]' 66,20
$arr[0,0] = 2
$arr[0,1] = "Comercializadora del Agro de Limarí"
$arr[0,2] = "Coquimbo"
$arr[0,3] = 44622
$arr[0,4] = 4
$arr[0,5] = "Fruta"
$arr[0,6] = 100103
$arr[0,7] = "Frutos de hueso (carozo)"
$arr[0,8] = 100103004
$arr[0,9] = "Durazno"
$arr[0,10] = "Doctor Davis"
$arr[0,11] = "Especial"
$arr[0,12] = 16
$arr[0,13] = 385000
$arr[0,14] = 390000
$arr[0,15] = 387500
$arr[0,16] = "`$/bins (400 kilos)"
$arr[0,17] = "Región de O'Higgins"
$arr[0,18] = 969
$arr[0,19] = 400
$arr[1,0] = 2
$arr[1,1] = "Comercializadora del Agro de Limarí"
$arr[1,2] = "Coquimbo"
$arr[1,3] = 44622
$arr[1,4] = 4
$arr[1,5] = "Fruta"
$arr[1,6] = 100103
$arr[1,7] = "Frutos de hueso (carozo)"
$arr[1,8] = 100103004
$arr[1,9] = "Durazno"
$arr[1,10] = "Doctor Davis"
$arr[1,11] = "Primera"
$arr[1,12] = 20
$arr[1,13] = 335000
$arr[1,14] = 340000
$arr[1,15] = 337500
$arr[1,16] = "`$/bins (400 kilos)"
$arr[1,17] = "Región de O'Higgins"
$arr[1,18] = 844
$arr[1,19] = 400
$arr[2,0] = 2
$arr[2,1] = "Comercializadora del Agro de Limarí"
$arr[2,2] = "Coquimbo"
$arr[2,3] = 44622
$arr[2,4] = 4
$arr[2,5] = "Fruta"
$arr[2,6] = 100103
$arr[2,7] = "Frutos de hueso (carozo)"
$arr[2,8] = 100103004
$arr[2,9] = "Durazno"
$arr[2,10] = "Doctor Davis"
$arr[2,11] = "Segunda"
$arr[2,12] = 20
$arr[2,13] = 295000
$arr[2,14] = 300000
$arr[2,15] = 297500
$arr[2,16] = "`$/bins (400 kilos)"
$arr[2,17] = "Región de O'Higgins"
$arr[2,18] = 744
$arr[2,19] = 400
$arr[3,0] = 2
$arr[3,1] = "Comercializadora del Agro de Limarí"
$arr[3,2] = "Coquimbo"
$arr[3,3] = 44203
$arr[3,4] = 4
$arr[3,5] = "Fruta"
$arr[3,6] = 100103
$arr[3,7] = "Frutos de hueso (carozo)"
$arr[3,8] = 100103004
$arr[3,9] = "Durazno"
$arr[3,10] = "Fortuna"
$arr[3,11] = "Especial"
$arr[3,12] = 200
$arr[3,13] = 19500
$arr[3,14] = 20000
$arr[3,15] = 19750
$arr[3,16] = "`$/caja 16 kilos empedrada"
$arr[3,17] = "Región Metropolitana"
$arr[3,18] = 1234
$arr[3,19] = 16
$arr[4,0] = 2
$arr[4,1] = "Comercializadora del Agro de Limarí"
$arr[4,2] = "Coquimbo"
$arr[4,3] = 44203
$arr[4,4] = 4
$arr[4,5] = "Fruta"
$arr[4,6] = 100103
$arr[4,7] = "Frutos de hueso (carozo)"
$arr[4,8] = 100103004
$arr[4,9] = "Durazno"
$arr[4,10] = "Fortuna"
$arr[4,11] = "Primera"
$arr[4,12] = 240
$arr[4,13] = 17500
$arr[4,14] = 18000
$arr[4,15] = 17750
$arr[4,16] = "`$/caja 16 kilos empedrada"
$arr[4,17] = "Región Metropolitana"
$arr[4,18] = 1109
$arr[4,19] = 16
$arr[5,0] = 2
$arr[5,1] = "Comercializadora del Agro de Limarí"
$arr[5,2] = "Coquimbo"
$arr[5,3] = 44203
$arr[5,4] = 4
$arr[5,5] = "Fruta"
$arr[5,6] = 100103
$arr[5,7] = "Frutos de hueso (carozo)"
$arr[5,8] = 100103004
$arr[5,9] = "Durazno"
$arr[5,10] = "Fortuna"
$arr[5,11] = "Segunda"
$arr[5,12] = 160
$arr[5,13] = 15500
$arr[5,14] = 16000
$arr[5,15] = 15750
$arr[5,16] = "`$/caja 16 kilos empedrada"
$arr[5,17] = "Región Metropolitana"
$arr[5,18] = 984
$arr[5,19] = 16
$arr[6,0] = 2
$arr[6,1] = "Comercializadora del Agro de Limarí"
$arr[6,2] = "Coquimbo"
$arr[6,3] = 44203
$arr[6,4] = 4
$arr[6,5] = "Fruta"
$arr[6,6] = 100103
$arr[6,7] = "Frutos de hueso (carozo)"
$arr[6,8] = 100103004
$arr[6,9] = "Durazno"
$arr[6,10] = "Toscana"
$arr[6,11] = "Especial"
$arr[6,12] = 200
$arr[6,13] = 19500
$arr[6,14] = 20000
$arr[6,15] = 19750
$arr[6,16] = "`$/caja 16 kilos empedrada"
$arr[6,17] = "Región de O'Higgins"
$arr[6,18] = 1234
$arr[6,19] = 16
$arr[7,0] = 2
$arr[7,1] = "Comercializadora del Agro de Limarí"
$arr[7,2] = "Coquimbo"
$arr[7,3] = 44203
$arr[7,4] = 4
$arr[7,5] = "Fruta"
$arr[7,6] = 100103
$arr[7,7] = "Frutos de hueso (carozo)"
$arr[7,8] = 100103004
$arr[7,9] = "Durazno"
$arr[7,10] = "Toscana"
$arr[7,11] = "Primera"
$arr[7,12] = 240
$arr[7,13] = 17500
$arr[7,14] = 18000
$arr[7,15] = 17750
$arr[7,16] = "`$/caja 16 kilos empedrada"
$arr[7,17] = "Región de O'Higgins"
$arr[7,18] = 1109
$arr[7,19] = 16
$arr[8,0] = 2
$arr[8,1] = "Comercializadora del Agro de Limarí"
$arr[8,2] = "Coquimbo"
$arr[8,3] = 44203
$arr[8,4] = 4
$arr[8,5] = "Fruta"
$arr[8,6] = 100103
$arr[8,7] = "Frutos de hueso (carozo)"
$arr[8,8] = 100103004
$arr[8,9] = "Durazno"
$arr[8,10] = "Toscana"
$arr[8,11] = "Segunda"
$arr[8,12] = 160
$arr[8,13] = 15500
$arr[8,14] = 16000
$arr[8,15] = 15750
$arr[8,16] = "`$/caja 16 kilos empedrada"
$arr[8,17] = "Región de O'Higgins"
$arr[8,18] = 984
$arr[8,19] = 16
$arr[9,0] = 2
$arr[9,1] = "Comercializadora del Agro de Limarí"
$arr[9,2] = "Coquimbo"
$arr[9,3] = 44245
$arr[9,4] = 4
$arr[9,5] = "Fruta"
$arr[9,6] = 100103
$arr[9,7] = "Frutos de hueso (carozo)"
$arr[9,8] = 100103004
$arr[9,9] = "Durazno"
$arr[9,10] = "Doctor Davis"
$arr[9,11] = "Especial"
$arr[9,12] = 140
$arr[9,13] = 17500
$arr[9,14] = 18000
$arr[9,15] = 17750
$arr[9,16] = "`$/caja 16 kilos empedrada"
$arr[9,17] = "Región de O'Higgins"
$arr[9,18] = 1109
$arr[9,19] = 16
$arr[10,0] = 2
$arr[10,1] = "Comercializadora del Agro de Limarí"
$arr[10,2] = "Coquimbo"
$arr[10,3] = 44245
$arr[10,4] = 4
$arr[10,5] = "Fruta"
$arr[10,6] = 100103
$arr[10,7] = "Frutos de hueso (carozo)"
$arr[10,8] = 100103004
$arr[10,9] = "Durazno"
$arr[10,10] = "Doctor Davis"
$arr[10,11] = "Primera"
$arr[10,12] = 200
$arr[10,13] = 15500
$arr[10,14] = 16000
$arr[10,15] = 15750
$arr[10,16] = "`$/caja 16 kilos empedrada"
$arr[10,17] = "Región de O'Higgins"
$arr[10,18] = 984
$arr[10,19] = 16
$arr[11,0] = 2
$arr[11,1] = "Comercializadora del Agro de Limarí"
$arr[11,2] = "Coquimbo"
$arr[11,3] = 44245
$arr[11,4] = 4
$arr[11,5] = "Fruta"
$arr[11,6] = 100103
$arr[11,7] = "Frutos de hueso (carozo)"
$arr[11,8] = 100103004
$arr[11,9] = "Durazno"
$arr[11,10] = "Doctor Davis"
$arr[11,11] = "Segunda"
$arr[11,12] = 160
$arr[11,13] = 13500
$arr[11,14] = 14000
$arr[11,15] = 13750
$arr[11,16] = "`$/caja 16 kilos empedrada"
$arr[11,17] = "Región de O'Higgins"
$arr[11,18] = 859
$arr[11,19] = 16
$arr[12,0] = 2
$arr[12,1] = "Comercializadora del Agro de Limarí"
$arr[12,2] = "Coquimbo"
$arr[12,3] = 44510
$arr[12,4] = 4
$arr[12,5] = "Fruta"
$arr[12,6] = 100103
$arr[12,7] = "Frutos de hueso (carozo)"
$arr[12,8] = 100103004
$arr[12,9] = "Durazno"
$arr[12,10] = "Florida King"
$arr[12,11] = "Primera"
$arr[12,12] = 360
$arr[12,13] = 12000
$arr[12,14] = 13000
$arr[12,15] = 12500
$arr[12,16] = "`$/bandeja 10 kilos granel"
$arr[12,17] = "Provincia de Limarí"
$arr[12,18] = 1250
$arr[12,19] = 10
$arr[13,0] = 2
$arr[13,1] = "Comercializadora del Agro de Limarí"
$arr[13,2] = "Coquimbo"
$arr[13,3] = 44510
$arr[13,4] = 4
$arr[13,5] = "Fruta"
$arr[13,6] = 100103
$arr[13,7] = "Frutos de hueso (carozo)"
$arr[13,8] = 100103004
$arr[13,9] = "Durazno"
$arr[13,10] = "Florida King"
$arr[13,11] = "Segunda"
$arr[13,12] = 260
$arr[13,13] = 9000
$arr[13,14] = 10000
$arr[13,15] = 9500
$arr[13,16] = "`$/bandeja 10 kilos granel"
$arr[13,17] = "Provincia de Limarí"
$arr[13,18] = 950
$arr[13,19] = 10
$arr[14,0] = 2
$arr[14,1] = "Comercializadora del Agro de Limarí"
$arr[14,2] = "Coquimbo"
$arr[14,3] = 44601
$arr[14,4] = 4
$arr[14,5] = "Fruta"
$arr[14,6] = 100103
$arr[14,7] = "Frutos de hueso (carozo)"
$arr[14,8] = 100103004
$arr[14,9] = "Durazno"
$arr[14,10] = "Loadel"
$arr[14,11] = "Especial"
$arr[14,12] = 100
$arr[14,13] = 22500
$arr[14,14] = 23000
$arr[14,15] = 22750
$arr[14,16] = "`$/caja 16 kilos empedrada"
$arr[14,17] = "Región de O'Higgins"
$arr[14,18] = 1422
$arr[14,19] = 16
$arr[15,0] = 2
$arr[15,1] = "Comercializadora del Agro de Limarí"
$arr[15,2] = "Coquimbo"
$arr[15,3] = 44601
$arr[15,4] = 4
$arr[15,5] = "Fruta"
$arr[15,6] = 100103
$arr[15,7] = "Frutos de hueso (carozo)"
$arr[15,8] = 100103004
$arr[15,9] = "Durazno"
$arr[15,10] = "Loadel"
$arr[15,11] = "Primera"
$arr[15,12] = 40
$arr[15,13] = 18500
$arr[15,14] = 19000
$arr[15,15] = 18750
$arr[15,16] = "`$/caja 16 kilos empedrada"
$arr[15,17] = "Región de O'Higgins"
$arr[15,18] = 1172
$arr[15,19] = 16
$arr[16,0] = 2
$arr[16,1] = "Comercializadora del Agro de Limarí"
$arr[16,2] = "Coquimbo"
$arr[16,3] = 44265
$arr[16,4] = 4
$arr[16,5] = "Fruta"
$arr[16,6] = 100103
$arr[16,7] = "Frutos de hueso (carozo)"
$arr[16,8] = 100103004
$arr[16,9] = "Durazno"
$arr[16,10] = "Phillips Cling"
$arr[16,11] = "Especial"
$arr[16,12] = 200
$arr[16,13] = 18500
$arr[16,14] = 19000
$arr[16,15] = 18750
$arr[16,16] = "`$/caja 16 kilos empedrada"
$arr[16,17] = "Región de O'Higgins"
$arr[16,18] = 1172
$arr[16,19] = 16
$arr[17,0] = 2
$arr[17,1] = "Comercializadora del Agro de Limarí"
$arr[17,2] = "Coquimbo"
$arr[17,3] = 44265
$arr[17,4] = 4
$arr[17,5] = "Fruta"
$arr[17,6] = 100103
$arr[17,7] = "Frutos de hueso (carozo)"
$arr[17,8] = 100103004
$arr[17,9] = "Durazno"
$arr[17,10] = "Phillips Cling"
$arr[17,11] = "Primera"
$arr[17,12] = 240
$arr[17,13] = 16500
$arr[17,14] = 17000
$arr[17,15] = 16750
$arr[17,16] = "`$/caja 16 kilos empedrada"
$arr[17,17] = "Región de O'Higgins"
$arr[17,18] = 1047
$arr[17,19] = 16
$arr[18,0] = 2
$arr[18,1] = "Comercializadora del Agro de Limarí"
$arr[18,2] = "Coquimbo"
$arr[18,3] = 44265
$arr[18,4] = 4
$arr[18,5] = "Fruta"
$arr[18,6] = 100103
$arr[18,7] = "Frutos de hueso (carozo)"
$arr[18,8] = 100103004
$arr[18,9] = "Durazno"
$arr[18,10] = "Phillips Cling"
$arr[18,11] = "Segunda"
$arr[18,12] = 300
$arr[18,13] = 13500
$arr[18,14] = 14000
$arr[18,15] = 13750
$arr[18,16] = "`$/caja 16 kilos empedrada"
$arr[18,17] = "Región de O'Higgins"
$arr[18,18] = 859
$arr[18,19] = 16
$arr[19,0] = 2
$arr[19,1] = "Comercializadora del Agro de Limarí"
$arr[19,2] = "Coquimbo"
$arr[19,3] = 44202
$arr[19,4] = 4
$arr[19,5] = "Fruta"
$arr[19,6] = 100103
$arr[19,7] = "Frutos de hueso (carozo)"
$arr[19,8] = 100103004
$arr[19,9] = "Durazno"
$arr[19,10] = "Andross"
$arr[19,11] = "Especial"
$arr[19,12] = 120
$arr[19,13] = 20000
$arr[19,14] = 21000
$arr[19,15] = 20500
$arr[19,16] = "`$/caja 16 kilos empedrada"
$arr[19,17] = "Región Metropolitana"
$arr[19,18] = 1281
$arr[19,19] = 16
$arr[20,0] = 2
$arr[20,1] = "Comercializadora del Agro de Limarí"
$arr[20,2] = "Coquimbo"
$arr[20,3] = 44202
$arr[20,4] = 4
$arr[20,5] = "Fruta"
$arr[20,6] = 100103
$arr[20,7] = "Frutos de hueso (carozo)"
$arr[20,8] = 100103004
$arr[20,9] = "Durazno"
$arr[20,10] = "Andross"
$arr[20,11] = "Primera"
$arr[20,12] = 200
$arr[20,13] = 18500
$arr[20,14] = 19000
$arr[20,15] = 18750
$arr[20,16] = "`$/caja 16 kilos empedrada"
$arr[20,17] = "Región Metropolitana"
$arr[20,18] = 1172
$arr[20,19] = 16
$arr[21,0] = 2
$arr[21,1] = "Comercializadora del Agro de Limarí"
$arr[21,2] = "Coquimbo"
$arr[21,3] = 44202
$arr[21,4] = 4
$arr[21,5] = "Fruta"
$arr[21,6] = 100103
$arr[21,7] = "Frutos de hueso (carozo)"
$arr[21,8] = 100103004
$arr[21,9] = "Durazno"
$arr[21,10] = "Andross"
$arr[21,11] = "Segunda"
$arr[21,12] = 200
$arr[21,13] = 17000
$arr[21,14] = 19500
$arr[21,15] = 18250
$arr[21,16] = "`$/caja 16 kilos empedrada"
$arr[21,17] = "Región Metropolitana"
$arr[21,18] = 1141
$arr[21,19] = 16
$arr[22,0] = 2
$arr[22,1] = "Comercializadora del Agro de Limarí"
$arr[22,2] = "Coquimbo"
$arr[22,3] = 44202
$arr[22,4] = 4
$arr[22,5] = "Fruta"
$arr[22,6] = 100103
$arr[22,7] = "Frutos de hueso (carozo)"
$arr[22,8] = 100103004
$arr[22,9] = "Durazno"
$arr[22,10] = "Fortuna"
$arr[22,11] = "Especial"
$arr[22,12] = 120
$arr[22,13] = 19500
$arr[22,14] = 20000
$arr[22,15] = 19750
$arr[22,16] = "`$/caja 16 kilos empedrada"
$arr[22,17] = "Región Metropolitana"
$arr[22,18] = 1234
$arr[22,19] = 16
$arr[23,0] = 2
$arr[23,1] = "Comercializadora del Agro de Limarí"
$arr[23,2] = "Coquimbo"
$arr[23,3] = 44202
$arr[23,4] = 4
$arr[23,5] = "Fruta"
$arr[23,6] = 100103
$arr[23,7] = "Frutos de hueso (carozo)"
$arr[23,8] = 100103004
$arr[23,9] = "Durazno"
$arr[23,10] = "Fortuna"
$arr[23,11] = "Primera"
$arr[23,12] = 200
$arr[23,13] = 17500
$arr[23,14] = 18000
$arr[23,15] = 17750
$arr[23,16] = "`$/caja 16 kilos empedrada"
$arr[23,17] = "Región Metropolitana"
$arr[23,18] = 1109
$arr[23,19] = 16
$arr[24,0] = 2
$arr[24,1] = "Comercializadora del Agro de Limarí"
$arr[24,2] = "Coquimbo"
$arr[24,3] = 44202
$arr[24,4] = 4
$arr[24,5] = "Fruta"
$arr[24,6] = 100103
$arr[24,7] = "Frutos de hueso (carozo)"
$arr[24,8] = 100103004
$arr[24,9] = "Durazno"
$arr[24,10] = "Fortuna"
$arr[24,11] = "Segunda"
$arr[24,12] = 200
$arr[24,13] = 15500
$arr[24,14] = 16000
$arr[24,15] = 15750
$arr[24,16] = "`$/caja 16 kilos empedrada"
$arr[24,17] = "Región Metropolitana"
$arr[24,18] = 984
$arr[24,19] = 16
$arr[25,0] = 2
$arr[25,1] = "Comercializadora del Agro de Limarí"
$arr[25,2] = "Coquimbo"
$arr[25,3] = 44210
$arr[25,4] = 4
$arr[25,5] = "Fruta"
$arr[25,6] = 100103
$arr[25,7] = "Frutos de hueso (carozo)"
$arr[25,8] = 100103004
$arr[25,9] = "Durazno"
$arr[25,10] = "Fortuna"
$arr[25,11] = "Especial"
$arr[25,12] = 160
$arr[25,13] = 17500
$arr[25,14] = 18000
$arr[25,15] = 17750
$arr[25,16] = "`$/caja 16 kilos empedrada"
$arr[25,17] = "Región Metropolitana"
$arr[25,18] = 1109
$arr[25,19] = 16
$arr[26,0] = 2
$arr[26,1] = "Comercializadora del Agro de Limarí"
$arr[26,2] = "Coquimbo"
$arr[26,3] = 44210
$arr[26,4] = 4
$arr[26,5] = "Fruta"
$arr[26,6] = 100103
$arr[26,7] = "Frutos de hueso (carozo)"
$arr[26,8] = 100103004
$arr[26,9] = "Durazno"
$arr[26,10] = "Fortuna"
$arr[26,11] = "Primera"
$arr[26,12] = 200
$arr[26,13] = 15500
$arr[26,14] = 16000
$arr[26,15] = 15750
$arr[26,16] = "`$/caja 16 kilos empedrada"
$arr[26,17] = "Región Metropolitana"
$arr[26,18] = 984
$arr[26,19] = 16
$arr[27,0] = 2
$arr[27,1] = "Comercializadora del Agro de Limarí"
$arr[27,2] = "Coquimbo"
$arr[27,3] = 44210
$arr[27,4] = 4
$arr[27,5] = "Fruta"
$arr[27,6] = 100103
$arr[27,7] = "Frutos de hueso (carozo)"
$arr[27,8] = 100103004
$arr[27,9] = "Durazno"
$arr[27,10] = "Fortuna"
$arr[27,11] = "Segunda"
$arr[27,12] = 200
$arr[27,13] = 12500
$arr[27,14] = 13000
$arr[27,15] = 12750
$arr[27,16] = "`$/caja 16 kilos empedrada"
$arr[27,17] = "Región Metropolitana"
$arr[27,18] = 797
$arr[27,19] = 16
$arr[28,0] = 2
$arr[28,1] = "Comercializadora del Agro de Limarí"
$arr[28,2] = "Coquimbo"
$arr[28,3] = 44217
$arr[28,4] = 4
$arr[28,5] = "Fruta"
$arr[28,6] = 100103
$arr[28,7] = "Frutos de hueso (carozo)"
$arr[28,8] = 100103004
$arr[28,9] = "Durazno"
$arr[28,10] = "Andross"
$arr[28,11] = "Especial"
$arr[28,12] = 240
$arr[28,13] = 19500
$arr[28,14] = 20000
$arr[28,15] = 19750
$arr[28,16] = "`$/caja 16 kilos empedrada"
$arr[28,17] = "Región de O'Higgins"
$arr[28,18] = 1234
$arr[28,19] = 16
$arr[29,0] = 2
$arr[29,1] = "Comercializadora del Agro de Limarí"
$arr[29,2] = "Coquimbo"
$arr[29,3] = 44217
$arr[29,4] = 4
$arr[29,5] = "Fruta"
$arr[29,6] = 100103
$arr[29,7] = "Frutos de hueso (carozo)"
$arr[29,8] = 100103004
$arr[29,9] = "Durazno"
$arr[29,10] = "Andross"
$arr[29,11] = "Primera"
$arr[29,12] = 240
$arr[29,13] = 17500
$arr[29,14] = 18000
$arr[29,15] = 17750
$arr[29,16] = "`$/caja 16 kilos empedrada"
$arr[29,17] = "Región de O'Higgins"
$arr[29,18] = 1109
$arr[29,19] = 16
$arr[30,0] = 2
$arr[30,1] = "Comercializadora del Agro de Limarí"
$arr[30,2] = "Coquimbo"
$arr[30,3] = 44217
$arr[30,4] = 4
$arr[30,5] = "Fruta"
$arr[30,6] = 100103
$arr[30,7] = "Frutos de hueso (carozo)"
$arr[30,8] = 100103004
$arr[30,9] = "Durazno"
$arr[30,10] = "Andross"
$arr[30,11] = "Segunda"
$arr[30,12] = 200
$arr[30,13] = 15500
$arr[30,14] = 16000
$arr[30,15] = 15750
$arr[30,16] = "`$/caja 16 kilos empedrada"
$arr[30,17] = "Región de O'Higgins"
$arr[30,18] = 984
$arr[30,19] = 16
$arr[31,0] = 2
$arr[31,1] = "Comercializadora del Agro de Limarí"
$arr[31,2] = "Coquimbo"
$arr[31,3] = 44238
$arr[31,4] = 4
$arr[31,5] = "Fruta"
$arr[31,6] = 100103
$arr[31,7] = "Frutos de hueso (carozo)"
$arr[31,8] = 100103004
$arr[31,9] = "Durazno"
$arr[31,10] = "Doctor Davis"
$arr[31,11] = "Especial"
$arr[31,12] = 130
$arr[31,13] = 18000
$arr[31,14] = 18500
$arr[31,15] = 18192
$arr[31,16] = "`$/caja 16 kilos empedrada"
$arr[31,17] = "Provincia de Limarí"
$arr[31,18] = 1137
$arr[31,19] = 16
$arr[32,0] = 2
$arr[32,1] = "Comercializadora del Agro de Limarí"
$arr[32,2] = "Coquimbo"
$arr[32,3] = 44238
$arr[32,4] = 4
$arr[32,5] = "Fruta"
$arr[32,6] = 100103
$arr[32,7] = "Frutos de hueso (carozo)"
$arr[32,8] = 100103004
$arr[32,9] = "Durazno"
$arr[32,10] = "Doctor Davis"
$arr[32,11] = "Primera"
$arr[32,12] = 160
$arr[32,13] = 16000
$arr[32,14] = 16500
$arr[32,15] = 16250
$arr[32,16] = "`$/caja 16 kilos empedrada"
$arr[32,17] = "Provincia de Limarí"
$arr[32,18] = 1016
$arr[32,19] = 16
$arr[33,0] = 2
$arr[33,1] = "Comercializadora del Agro de Limarí"
$arr[33,2] = "Coquimbo"
$arr[33,3] = 44238
$arr[33,4] = 4
$arr[33,5] = "Fruta"
$arr[33,6] = 100103
$arr[33,7] = "Frutos de hueso (carozo)"
$arr[33,8] = 100103004
$arr[33,9] = "Durazno"
$arr[33,10] = "Doctor Davis"
$arr[33,11] = "Segunda"
$arr[33,12] = 120
$arr[33,13] = 14000
$arr[33,14] = 14500
$arr[33,15] = 14250
$arr[33,16] = "`$/caja 16 kilos empedrada"
$arr[33,17] = "Provincia de Limarí"
$arr[33,18] = 891
$arr[33,19] = 16
$arr[34,0] = 2
$arr[34,1] = "Comercializadora del Agro de Limarí"
$arr[34,2] = "Coquimbo"
$arr[34,3] = 44175
$arr[34,4] = 4
$arr[34,5] = "Fruta"
$arr[34,6] = 100103
$arr[34,7] = "Frutos de hueso (carozo)"
$arr[34,8] = 100103004
$arr[34,9] = "Durazno"
$arr[34,10] = "Flavor Crest"
$arr[34,11] = "Primera"
$arr[34,12] = 200
$arr[34,13] = 17500
$arr[34,14] = 18000
$arr[34,15] = 17750
$arr[34,16] = "`$/caja 16 kilos empedrada"
$arr[34,17] = "Región de O'Higgins"
$arr[34,18] = 1109
$arr[34,19] = 16
$arr[35,0] = 2
$arr[35,1] = "Comercializadora del Agro de Limarí"
$arr[35,2] = "Coquimbo"
$arr[35,3] = 44175
$arr[35,4] = 4
$arr[35,5] = "Fruta"
$arr[35,6] = 100103
$arr[35,7] = "Frutos de hueso (carozo)"
$arr[35,8] = 100103004
$arr[35,9] = "Durazno"
$arr[35,10] = "Flavor Crest"
$arr[35,11] = "Segunda"
$arr[35,12] = 240
$arr[35,13] = 14500
$arr[35,14] = 15000
$arr[35,15] = 14750
$arr[35,16] = "`$/caja 16 kilos empedrada"
$arr[35,17] = "Región de O'Higgins"
$arr[35,18] = 922
$arr[35,19] = 16
$arr[36,0] = 2
$arr[36,1] = "Comercializadora del Agro de Limarí"
$arr[36,2] = "Coquimbo"
$arr[36,3] = 44237
$arr[36,4] = 4
$arr[36,5] = "Fruta"
$arr[36,6] = 100103
$arr[36,7] = "Frutos de hueso (carozo)"
$arr[36,8] = 100103004
$arr[36,9] = "Durazno"
$arr[36,10] = "Doctor Davis"
$arr[36,11] = "Especial"
$arr[36,12] = 100
$arr[36,13] = 18500
$arr[36,14] = 19000
$arr[36,15] = 18750
$arr[36,16] = "`$/caja 16 kilos empedrada"
$arr[36,17] = "Región Metropolitana"
$arr[36,18] = 1172
$arr[36,19] = 16
$arr[37,0] = 2
$arr[37,1] = "Comercializadora del Agro de Limarí"
$arr[37,2] = "Coquimbo"
$arr[37,3] = 44237
$arr[37,4] = 4
$arr[37,5] = "Fruta"
$arr[37,6] = 100103
$arr[37,7] = "Frutos de hueso (carozo)"
$arr[37,8] = 100103004
$arr[37,9] = "Durazno"
$arr[37,10] = "Doctor Davis"
$arr[37,11] = "Primera"
$arr[37,12] = 140
$arr[37,13] = 16500
$arr[37,14] = 17000
$arr[37,15] = 16750
$arr[37,16] = "`$/caja 16 kilos empedrada"
$arr[37,17] = "Región Metropolitana"
$arr[37,18] = 1047
$arr[37,19] = 16
$arr[38,0] = 2
$arr[38,1] = "Comercializadora del Agro de Limarí"
$arr[38,2] = "Coquimbo"
$arr[38,3] = 44237
$arr[38,4] = 4
$arr[38,5] = "Fruta"
$arr[38,6] = 100103
$arr[38,7] = "Frutos de hueso (carozo)"
$arr[38,8] = 100103004
$arr[38,9] = "Durazno"
$arr[38,10] = "Doctor Davis"
$arr[38,11] = "Segunda"
$arr[38,12] = 200
$arr[38,13] = 14500
$arr[38,14] = 15000
$arr[38,15] = 14750
$arr[38,16] = "`$/caja 16 kilos empedrada"
$arr[38,17] = "Región Metropolitana"
$arr[38,18] = 922
$arr[38,19] = 16
$arr[39,0] = 2
$arr[39,1] = "Comercializadora del Agro de Limarí"
$arr[39,2] = "Coquimbo"
$arr[39,3] = 44609
$arr[39,4] = 4
$arr[39,5] = "Fruta"
$arr[39,6] = 100103
$arr[39,7] = "Frutos de hueso (carozo)"
$arr[39,8] = 100103004
$arr[39,9] = "Durazno"
$arr[39,10] = "Carson"
$arr[39,11] = "Especial"
$arr[39,12] = 200
$arr[39,13] = 17000
$arr[39,14] = 18000
$arr[39,15] = 17500
$arr[39,16] = "`$/caja 16 kilos empedrada"
$arr[39,17] = "Región de O'Higgins"
$arr[39,18] = 1094
$arr[39,19] = 16
$arr[40,0] = 2
$arr[40,1] = "Comercializadora del Agro de Limarí"
$arr[40,2] = "Coquimbo"
$arr[40,3] = 44609
$arr[40,4] = 4
$arr[40,5] = "Fruta"
$arr[40,6] = 100103
$arr[40,7] = "Frutos de hueso (carozo)"
$arr[40,8] = 100103004
$arr[40,9] = "Durazno"
$arr[40,10] = "Carson"
$arr[40,11] = "Primera"
$arr[40,12] = 200
$arr[40,13] = 15000
$arr[40,14] = 16000
$arr[40,15] = 15500
$arr[40,16] = "`$/caja 16 kilos empedrada"
$arr[40,17] = "Región de O'Higgins"
$arr[40,18] = 969
$arr[40,19] = 16
$arr[41,0] = 2
$arr[41,1] = "Comercializadora del Agro de Limarí"
$arr[41,2] = "Coquimbo"
$arr[41,3] = 44609
$arr[41,4] = 4
$arr[41,5] = "Fruta"
$arr[41,6] = 100103
$arr[41,7] = "Frutos de hueso (carozo)"
$arr[41,8] = 100103004
$arr[41,9] = "Durazno"
$arr[41,10] = "Carson"
$arr[41,11] = "Segunda"
$arr[41,12] = 200
$arr[41,13] = 13000
$arr[41,14] = 14000
$arr[41,15] = 13500
$arr[41,16] = "`$/caja 16 kilos empedrada"
$arr[41,17] = "Región de O'Higgins"
$arr[41,18] = 844
$arr[41,19] = 16
$arr[42,0] = 2
$arr[42,1] = "Comercializadora del Agro de Limarí"
$arr[42,2] = "Coquimbo"
$arr[42,3] = 44195
$arr[42,4] = 4
$arr[42,5] = "Fruta"
$arr[42,6] = 100103
$arr[42,7] = "Frutos de hueso (carozo)"
$arr[42,8] = 100103004
$arr[42,9] = "Durazno"
$arr[42,10] = "Flavor Crest"
$arr[42,11] = "Especial"
$arr[42,12] = 360
$arr[42,13] = 17500
$arr[42,14] = 18000
$arr[42,15] = 17750
$arr[42,16] = "`$/caja 16 kilos empedrada"
$arr[42,17] = "Región Metropolitana"
$arr[42,18] = 1109
$arr[42,19] = 16
$arr[43,0] = 2
$arr[43,1] = "Comercializadora del Agro de Limarí"
$arr[43,2] = "Coquimbo"
$arr[43,3] = 44195
$arr[43,4] = 4
$arr[43,5] = "Fruta"
$arr[43,6] = 100103
$arr[43,7] = "Frutos de hueso (carozo)"
$arr[43,8] = 100103004
$arr[43,9] = "Durazno"
$arr[43,10] = "Flavor Crest"
$arr[43,11] = "Primera"
$arr[43,12] = 300
$arr[43,13] = 15500
$arr[43,14] = 16000
$arr[43,15] = 15750
$arr[43,16] = "`$/caja 16 kilos empedrada"
$arr[43,17] = "Región Metropolitana"
$arr[43,18] = 984
$arr[43,19] = 16
$arr[44,0] = 2
$arr[44,1] = "Comercializadora del Agro de Limarí"
$arr[44,2] = "Coquimbo"
$arr[44,3] = 44195
$arr[44,4] = 4
$arr[44,5] = "Fruta"
$arr[44,6] = 100103
$arr[44,7] = "Frutos de hueso (carozo)"
$arr[44,8] = 100103004
$arr[44,9] = "Durazno"
$arr[44,10] = "Flavor Crest"
$arr[44,11] = "Segunda"
$arr[44,12] = 280
$arr[44,13] = 13500
$arr[44,14] = 14000
$arr[44,15] = 13750
$arr[44,16] = "`$/caja 16 kilos empedrada"
$arr[44,17] = "Región Metropolitana"
$arr[44,18] = 859
$arr[44,19] = 16
$arr[45,0] = 2
$arr[45,1] = "Comercializadora del Agro de Limarí"
$arr[45,2] = "Coquimbo"
$arr[45,3] = 44195
$arr[45,4] = 4
$arr[45,5] = "Fruta"
$arr[45,6] = 100103
$arr[45,7] = "Frutos de hueso (carozo)"
$arr[45,8] = 100103004
$arr[45,9] = "Durazno"
$arr[45,10] = "Royal Glory"
$arr[45,11] = "Especial"
$arr[45,12] = 340
$arr[45,13] = 17500
$arr[45,14] = 18000
$arr[45,15] = 17750
$arr[45,16] = "`$/caja 16 kilos empedrada"
$arr[45,17] = "Región Metropolitana"
$arr[45,18] = 1109
$arr[45,19] = 16
$arr[46,0] = 2
$arr[46,1] = "Comercializadora del Agro de Limarí"
$arr[46,2] = "Coquimbo"
$arr[46,3] = 44195
$arr[46,4] = 4
$arr[46,5] = "Fruta"
$arr[46,6] = 100103
$arr[46,7] = "Frutos de hueso (carozo)"
$arr[46,8] = 100103004
$arr[46,9] = "Durazno"
$arr[46,10] = "Royal Glory"
$arr[46,11] = "Primera"
$arr[46,12] = 240
$arr[46,13] = 15500
$arr[46,14] = 16000
$arr[46,15] = 15750
$arr[46,16] = "`$/caja 16 kilos empedrada"
$arr[46,17] = "Región Metropolitana"
$arr[46,18] = 984
$arr[46,19] = 16
$arr[47,0] = 2
$arr[47,1] = "Comercializadora del Agro de Limarí"
$arr[47,2] = "Coquimbo"
$arr[47,3] = 44195
$arr[47,4] = 4
$arr[47,5] = "Fruta"
$arr[47,6] = 100103
$arr[47,7] = "Frutos de hueso (carozo)"
$arr[47,8] = 100103004
$arr[47,9] = "Durazno"
$arr[47,10] = "Royal Glory"
$arr[47,11] = "Segunda"
$arr[47,12] = 220
$arr[47,13] = 13500
$arr[47,14] = 14000
$arr[47,15] = 13750
$arr[47,16] = "`$/caja 16 kilos empedrada"
$arr[47,17] = "Región Metropolitana"
$arr[47,18] = 859
$arr[47,19] = 16
$arr[48,0] = 2
$arr[48,1] = "Comercializadora del Agro de Limarí"
$arr[48,2] = "Coquimbo"
$arr[48,3] = 44252
$arr[48,4] = 4
$arr[48,5] = "Fruta"
$arr[48,6] = 100103
$arr[48,7] = "Frutos de hueso (carozo)"
$arr[48,8] = 100103004
$arr[48,9] = "Durazno"
$arr[48,10] = "Kakamas"
$arr[48,11] = "Especial"
$arr[48,12] = 240
$arr[48,13] = 20500
$arr[48,14] = 21000
$arr[48,15] = 20750
$arr[48,16] = "`$/caja 18 kilos empedrada"
$arr[48,17] = "Región de O'Higgins"
$arr[48,18] = 1153
$arr[48,19] = 18
$arr[49,0] = 2
$arr[49,1] = "Comercializadora del Agro de Limarí"
$arr[49,2] = "Coquimbo"
$arr[49,3] = 44252
$arr[49,4] = 4
$arr[49,5] = "Fruta"
$arr[49,6] = 100103
$arr[49,7] = "Frutos de hueso (carozo)"
$arr[49,8] = 100103004
$arr[49,9] = "Durazno"
$arr[49,10] = "Kakamas"
$arr[49,11] = "Primera"
$arr[49,12] = 240
$arr[49,13] = 17500
$arr[49,14] = 18000
$arr[49,15] = 17750
$arr[49,16] = "`$/caja 18 kilos empedrada"
$arr[49,17] = "Región de O'Higgins"
$arr[49,18] = 986
$arr[49,19] = 18
$arr[50,0] = 2
$arr[50,1] = "Comercializadora del Agro de Limarí"
$arr[50,2] = "Coquimbo"
$arr[50,3] = 44252
$arr[50,4] = 4
$arr[50,5] = "Fruta"
$arr[50,6] = 100103
$arr[50,7] = "Frutos de hueso (carozo)"
$arr[50,8] = 100103004
$arr[50,9] = "Durazno"
$arr[50,10] = "Kakamas"
$arr[50,11] = "Segunda"
$arr[50,12] = 200
$arr[50,13] = 14500
$arr[50,14] = 15000
$arr[50,15] = 14750
$arr[50,16] = "`$/caja 18 kilos empedrada"
$arr[50,17] = "Región de O'Higgins"
$arr[50,18] = 819
$arr[50,19] = 18
$arr[51,0] = 2
$arr[51,1] = "Comercializadora del Agro de Limarí"
$arr[51,2] = "Coquimbo"
$arr[51,3] = 44231
$arr[51,4] = 4
$arr[51,5] = "Fruta"
$arr[51,6] = 100103
$arr[51,7] = "Frutos de hueso (carozo)"
$arr[51,8] = 100103004
$arr[51,9] = "Durazno"
$arr[51,10] = "Doctor Davis"
$arr[51,11] = "Especial"
$arr[51,12] = 100
$arr[51,13] = 16500
$arr[51,14] = 17000
$arr[51,15] = 16750
$arr[51,16] = "`$/caja 16 kilos empedrada"
$arr[51,17] = "Región Metropolitana"
$arr[51,18] = 1047
$arr[51,19] = 16
$arr[52,0] = 2
$arr[52,1] = "Comercializadora del Agro de Limarí"
$arr[52,2] = "Coquimbo"
$arr[52,3] = 44231
$arr[52,4] = 4
$arr[52,5] = "Fruta"
$arr[52,6] = 100103
$arr[52,7] = "Frutos de hueso (carozo)"
$arr[52,8] = 100103004
$arr[52,9] = "Durazno"
$arr[52,10] = "Doctor Davis"
$arr[52,11] = "Primera"
$arr[52,12] = 120
$arr[52,13] = 14500
$arr[52,14] = 15000
$arr[52,15] = 14750
$arr[52,16] = "`$/caja 16 kilos empedrada"
$arr[52,17] = "Región Metropolitana"
$arr[52,18] = 922
$arr[52,19] = 16
$arr[53,0] = 2
$arr[53,1] = "Comercializadora del Agro de Limarí"
$arr[53,2] = "Coquimbo"
$arr[53,3] = 44231
$arr[53,4] = 4
$arr[53,5] = "Fruta"
$arr[53,6] = 100103
$arr[53,7] = "Frutos de hueso (carozo)"
$arr[53,8] = 100103004
$arr[53,9] = "Durazno"
$arr[53,10] = "Doctor Davis"
$arr[53,11] = "Segunda"
$arr[53,12] = 120
$arr[53,13] = 12500
$arr[53,14] = 13000
$arr[53,15] = 12750
$arr[53,16] = "`$/caja 16 kilos empedrada"
$arr[53,17] = "Región Metropolitana"
$arr[53,18] = 797
$arr[53,19] = 16
$arr[54,0] = 2
$arr[54,1] = "Comercializadora del Agro de Limarí"
$arr[54,2] = "Coquimbo"
$arr[54,3] = 44615
$arr[54,4] = 4
$arr[54,5] = "Fruta"
$arr[54,6] = 100103
$arr[54,7] = "Frutos de hueso (carozo)"
$arr[54,8] = 100103004
$arr[54,9] = "Durazno"
$arr[54,10] = "Doctor Davis"
$arr[54,11] = "Especial"
$arr[54,12] = 16
$arr[54,13] = 355000
$arr[54,14] = 360000
$arr[54,15] = 357500
$arr[54,16] = "`$/bins (400 kilos)"
$arr[54,17] = "Región de O'Higgins"
$arr[54,18] = 894
$arr[54,19] = 400
$arr[55,0] = 2
$arr[55,1] = "Comercializadora del Agro de Limarí"
$arr[55,2] = "Coquimbo"
$arr[55,3] = 44615
$arr[55,4] = 4
$arr[55,5] = "Fruta"
$arr[55,6] = 100103
$arr[55,7] = "Frutos de hueso (carozo)"
$arr[55,8] = 100103004
$arr[55,9] = "Durazno"
$arr[55,10] = "Doctor Davis"
$arr[55,11] = "Primera"
$arr[55,12] = 20
$arr[55,13] = 315000
$arr[55,14] = 320000
$arr[55,15] = 317500
$arr[55,16] = "`$/bins (400 kilos)"
$arr[55,17] = "Región de O'Higgins"
$arr[55,18] = 794
$arr[55,19] = 400
$arr[56,0] = 2
$arr[56,1] = "Comercializadora del Agro de Limarí"
$arr[56,2] = "Coquimbo"
$arr[56,3] = 44615
$arr[56,4] = 4
$arr[56,5] = "Fruta"
$arr[56,6] = 100103
$arr[56,7] = "Frutos de hueso (carozo)"
$arr[56,8] = 100103004
$arr[56,9] = "Durazno"
$arr[56,10] = "September Snow"
$arr[56,11] = "Especial"
$arr[56,12] = 10
$arr[56,13] = 325000
$arr[56,14] = 330000
$arr[56,15] = 327500
$arr[56,16] = "`$/bins (400 kilos)"
$arr[56,17] = "Región de O'Higgins"
$arr[56,18] = 819
$arr[56,19] = 400
$arr[57,0] = 2
$arr[57,1] = "Comercializadora del Agro de Limarí"
$arr[57,2] = "Coquimbo"
$arr[57,3] = 44615
$arr[57,4] = 4
$arr[57,5] = "Fruta"
$arr[57,6] = 100103
$arr[57,7] = "Frutos de hueso (carozo)"
$arr[57,8] = 100103004
$arr[57,9] = "Durazno"
$arr[57,10] = "September Snow"
$arr[57,11] = "Primera"
$arr[57,12] = 20
$arr[57,13] = 295000
$arr[57,14] = 300000
$arr[57,15] = 297500
$arr[57,16] = "`$/bins (400 kilos)"
$arr[57,17] = "Región de O'Higgins"
$arr[57,18] = 744
$arr[57,19] = 400
$arr[58,0] = 2
$arr[58,1] = "Comercializadora del Agro de Limarí"
$arr[58,2] = "Coquimbo"
$arr[58,3] = 44167
$arr[58,4] = 4
$arr[58,5] = "Fruta"
$arr[58,6] = 100103
$arr[58,7] = "Frutos de hueso (carozo)"
$arr[58,8] = 100103004
$arr[58,9] = "Durazno"
$arr[58,10] = "Springcrest"
$arr[58,11] = "Especial"
$arr[58,12] = 100
$arr[58,13] = 17500
$arr[58,14] = 18000
$arr[58,15] = 17750
$arr[58,16] = "`$/caja 16 kilos empedrada"
$arr[58,17] = "Región Metropolitana"
$arr[58,18] = 1109
$arr[58,19] = 16
$arr[59,0] = 2
$arr[59,1] = "Comercializadora del Agro de Limarí"
$arr[59,2] = "Coquimbo"
$arr[59,3] = 44167
$arr[59,4] = 4
$arr[59,5] = "Fruta"
$arr[59,6] = 100103
$arr[59,7] = "Frutos de hueso (carozo)"
$arr[59,8] = 100103004
$arr[59,9] = "Durazno"
$arr[59,10] = "Springcrest"
$arr[59,11] = "Primera"
$arr[59,12] = 200
$arr[59,13] = 15500
$arr[59,14] = 16000
$arr[59,15] = 15750
$arr[59,16] = "`$/caja 16 kilos empedrada"
$arr[59,17] = "Región Metropolitana"
$arr[59,18] = 984
$arr[59,19] = 16
$arr[60,0] = 2
$arr[60,1] = "Comercializadora del Agro de Limarí"
$arr[60,2] = "Coquimbo"
$arr[60,3] = 44167
$arr[60,4] = 4
$arr[60,5] = "Fruta"
$arr[60,6] = 100103
$arr[60,7] = "Frutos de hueso (carozo)"
$arr[60,8] = 100103004
$arr[60,9] = "Durazno"
$arr[60,10] = "Springcrest"
$arr[60,11] = "Segunda"
$arr[60,12] = 120
$arr[60,13] = 13500
$arr[60,14] = 14000
$arr[60,15] = 13750
$arr[60,16] = "`$/caja 16 kilos empedrada"
$arr[60,17] = "Región Metropolitana"
$arr[60,18] = 859
$arr[60,19] = 16
$arr[61,0] = 2
$arr[61,1] = "Comercializadora del Agro de Limarí"
$arr[61,2] = "Coquimbo"
$arr[61,3] = 44258
$arr[61,4] = 4
$arr[61,5] = "Fruta"
$arr[61,6] = 100103
$arr[61,7] = "Frutos de hueso (carozo)"
$arr[61,8] = 100103004
$arr[61,9] = "Durazno"
$arr[61,10] = "Doctor Davis"
$arr[61,11] = "Especial"
$arr[61,12] = 240
$arr[61,13] = 18500
$arr[61,14] = 19000
$arr[61,15] = 18750
$arr[61,16] = "`$/caja 16 kilos empedrada"
$arr[61,17] = "Región de O'Higgins"
$arr[61,18] = 1172
$arr[61,19] = 16
$arr[62,0] = 2
$arr[62,1] = "Comercializadora del Agro de Limarí"
$arr[62,2] = "Coquimbo"
$arr[62,3] = 44258
$arr[62,4] = 4
$arr[62,5] = "Fruta"
$arr[62,6] = 100103
$arr[62,7] = "Frutos de hueso (carozo)"
$arr[62,8] = 100103004
$arr[62,9] = "Durazno"
$arr[62,10] = "Doctor Davis"
$arr[62,11] = "Primera"
$arr[62,12] = 240
$arr[62,13] = 16500
$arr[62,14] = 17000
$arr[62,15] = 16750
$arr[62,16] = "`$/caja 16 kilos empedrada"
$arr[62,17] = "Región de O'Higgins"
$arr[62,18] = 1047
$arr[62,19] = 16
$arr[63,0] = 2
$arr[63,1] = "Comercializadora del Agro de Limarí"
$arr[63,2] = "Coquimbo"
$arr[63,3] = 44595
$arr[63,4] = 4
$arr[63,5] = "Fruta"
$arr[63,6] = 100103
$arr[63,7] = "Frutos de hueso (carozo)"
$arr[63,8] = 100103004
$arr[63,9] = "Durazno"
$arr[63,10] = "Andross"
$arr[63,11] = "Especial"
$arr[63,12] = 200
$arr[63,13] = 25000
$arr[63,14] = 26000
$arr[63,15] = 25500
$arr[63,16] = "`$/caja 18 kilos empedrada"
$arr[63,17] = "Región de O'Higgins"
$arr[63,18] = 1417
$arr[63,19] = 18
$arr[64,0] = 2
$arr[64,1] = "Comercializadora del Agro de Limarí"
$arr[64,2] = "Coquimbo"
$arr[64,3] = 44595
$arr[64,4] = 4
$arr[64,5] = "Fruta"
$arr[64,6] = 100103
$arr[64,7] = "Frutos de hueso (carozo)"
$arr[64,8] = 100103004
$arr[64,9] = "Durazno"
$arr[64,10] = "Andross"
$arr[64,11] = "Primera"
$arr[64,12] = 240
$arr[64,13] = 20000
$arr[64,14] = 21000
$arr[64,15] = 20500
$arr[64,16] = "`$/caja 18 kilos empedrada"
$arr[64,17] = "Región de O'Higgins"
$arr[64,18] = 1139
$arr[64,19] = 18
$arr[65,0] = 2
$arr[65,1] = "Comercializadora del Agro de Limarí"
$arr[65,2] = "Coquimbo"
$arr[65,3] = 44595
$arr[65,4] = 4
$arr[65,5] = "Fruta"
$arr[65,6] = 100103
$arr[65,7] = "Frutos de hueso (carozo)"
$arr[65,8] = 100103004
$arr[65,9] = "Durazno"
$arr[65,10] = "Andross"
$arr[65,11] = "Segunda"
$arr[65,12] = 300
$arr[65,13] = 15000
$arr[65,14] = 16000
$arr[65,15] = 15500
$arr[65,16] = "`$/caja 18 kilos empedrada"
$arr[65,17] = "Región de O'Higgins"
$arr[65,18] = 861
$arr[65,19] = 18

$ws.Range("A45:T110").Value = $arr

# Ensure the date column (D) keeps the existing date number format for all rows, including newly added 108:110
$ws.Range("D45:D110").NumberFormat = $ws.Range("D44").NumberFormat
